# Apply cryptos list update (prices and 1h volume %) per commit "Updated cryptos list on Thu Nov 14 17:30:39 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "88.785.77"
$ws.Range("E2").Value = "  -3.88%  "

$ws.Range("D3").Value = "3.135.70"
$ws.Range("E3").Value = "  -4.43%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'214.79"
$ws.Range("E5").Value = "  -1.31%  "

$ws.Range("D6").Value = "'633.58"
$ws.Range("E6").Value = "  +0.53%  "

$ws.Range("D7").Value = "'0.393"
$ws.Range("E7").Value = "  -4.67%  "

$ws.Range("D8").Value = "'0.750"
$ws.Range("E8").Value = "  +4.25%  "

$ws.Range("E9").Value = "  +0.13%  "

$ws.Range("D10").Value = "3.132.29"
$ws.Range("E10").Value = "  -4.38%  "

$ws.Range("D11").Value = "'0.557"
$ws.Range("E11").Value = "  -5.23%  "

$ws.Range("E12").Value = "  -0.62%  "

$ws.Range("D13").Value = "'0.0000251"
$ws.Range("E13").Value = "  -4.90%  "

$ws.Range("D14").Value = "'5.31"
$ws.Range("E14").Value = "  -0.47%  "

$ws.Range("D15").Value = "88.731.97"
$ws.Range("E15").Value = "  -3.73%  "

$ws.Range("D16").Value = "3.722.65"
$ws.Range("E16").Value = "  -4.13%  "

$ws.Range("D17").Value = "'32.39"
$ws.Range("E17").Value = "  -5.36%  "

$ws.Range("D18").Value = "3.151.93"
$ws.Range("E18").Value = "  -4.79%  "

# PEPE/SuiNetwork swapped ranking positions
$ws.Range("B19").Value = "SuiNetwork"
$ws.Range("C19").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D19").Value = "'3.39"
$ws.Range("E19").Value = "  +1.61%  "

$ws.Range("B20").Value = "PEPE"
$ws.Range("C20").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D20").Value = "'0.0000228"
$ws.Range("E20").Value = "  +19.34%  "

$ws.Range("D21").Value = "'13.30"
$ws.Range("E21").Value = "  -5.19%  "

$ws.Range("D22").Value = "'426.09"
$ws.Range("E22").Value = "  -3.31%  "

$ws.Range("D23").Value = "'8.38"
$ws.Range("E23").Value = "  -5.69%  "

$ws.Range("D24").Value = "'4.89"
$ws.Range("E24").Value = "  -7.29%  "

$ws.Range("D25").Value = "'5.42"
$ws.Range("E25").Value = "  +0.56%  "

$ws.Range("D26").Value = "'11.54"
$ws.Range("E26").Value = "  -5.85%  "

$ws.Range("D27").Value = "'81.54"
$ws.Range("E27").Value = "  +5.94%  "

$ws.Range("D28").Value = "3.311.33"
$ws.Range("E28").Value = "  -4.64%  "

$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.13%  "

$ws.Range("D30").Value = "'0.158"
$ws.Range("E30").Value = "  -11.90%  "

$ws.Range("D31").Value = "'0.996"
$ws.Range("E31").Value = "  -0.68%  "

$ws.Range("D32").Value = "'4.04"
$ws.Range("E32").Value = "  +11.61%  "

$ws.Range("D33").Value = "'8.18"
$ws.Range("E33").Value = "  -6.63%  "

$ws.Range("D34").Value = "'510.68"
$ws.Range("E34").Value = "  -8.08%  "

$ws.Range("D35").Value = "'7.10"
$ws.Range("E35").Value = "  -0.84%  "

# Fetch.AI/Kaspa swapped ranking positions
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "'0.147"
$ws.Range("E36").Value = "  +12.09%  "

$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").Value = "'1.30"
$ws.Range("E37").Value = "  +1.67%  "

$ws.Range("E38").Value = "  -4.80%  "

$ws.Range("D39").Value = "'21.90"
$ws.Range("E39").Value = "  -3.44%  "

$ws.Range("D40").Value = "'22.24"
$ws.Range("E40").Value = "  -1.14%  "

$ws.Range("E41").Value = "  +0.29%  "

$ws.Range("E42").Value = "  -0.07%  "

$ws.Range("D43").Value = "'1.86"
$ws.Range("E43").Value = "  -7.01%  "

$ws.Range("D44").Value = "'0.364"
$ws.Range("E44").Value = "  -7.23%  "

$ws.Range("D45").Value = "'145.89"
$ws.Range("E45").Value = "  -3.20%  "

$ws.Range("D46").Value = "'43.75"
$ws.Range("E46").Value = "  -3.72%  "

$ws.Range("D47").Value = "'0.128"
$ws.Range("E47").Value = "  -2.09%  "

$ws.Range("D48").Value = "'164.98"
$ws.Range("E48").Value = "  -8.54%  "

$ws.Range("D49").Value = "'0.720"
$ws.Range("E49").Value = "  -1.34%  "

$ws.Range("D50").Value = "'24.40"
$ws.Range("E50").Value = "  -2.37%  "

$ws.Range("D51").Value = "'0.598"
$ws.Range("E51").Value = "  -5.58%  "
